$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.3193623891079
$ws.Range("C2").Value = 12.4024333816592
$ws.Range("D2").Value = 14.2362913965566
$ws.Range("C3").Value = 10.4568787545218
$ws.Range("D3").Value = 13.1205295508446
$ws.Range("C4").Value = 12.3421846754442
$ws.Range("D4").Value = 19.3723269050858
$ws.Range("C5").Value = 9.32709677839388
$ws.Range("D5").Value = 15.0300437350976
$ws.Range("C6").Value = 6.96753761108763
$ws.Range("D6").Value = 13.3983121510222
$ws.Range("C7").Value = 4.66208118698497
$ws.Range("D7").Value = 11.6779331570705
$ws.Range("C8").Value = 5.59256483368605
$ws.Range("D8").Value = 12.0052939618059
$ws.Range("C9").Value = 8.98625587485825
$ws.Range("D9").Value = 15.4399363649316
$ws.Range("C10").Value = 6.48114033121892
$ws.Range("D10").Value = 25.4187079716149
$ws.Range("B11").Value = 16.9764659249636
$ws.Range("C11").Value = 14.3625021905328
$ws.Range("D11").Value = 19.5904296593944
$ws.Range("B12").Value = 2.59148705771949
$ws.Range("C12").Value = 0.268300024509501
$ws.Range("D12").Value = 4.91467409092947
$ws.Range("B13").Value = 14.5287976377778
$ws.Range("C13").Value = 13.0027964933064
$ws.Range("D13").Value = 16.0547987822492
$ws.Range("B14").Value = 5.28551547481145
$ws.Range("C14").Value = 4.68172699828632
$ws.Range("D14").Value = 5.88930395133658
$ws.Range("C15").Value = 2.68306425428489
$ws.Range("D15").Value = 4.18749594129651
$ws.Range("C16").Value = 2.69481326877843
$ws.Range("D16").Value = 6.78427766287268
$ws.Range("C17").Value = 1.71622282925407
$ws.Range("D17").Value = 4.81577102495517
$ws.Range("C18").Value = 1.46382663299961
$ws.Range("D18").Value = 5.31047236203405
$ws.Range("C19").Value = 0.809705039590391
$ws.Range("D19").Value = 5.17302134030723
$ws.Range("C20").Value = 0.712075912213744
$ws.Range("D20").Value = 4.22439972061379
$ws.Range("C21").Value = 1.35671177269789
$ws.Range("D21").Value = 4.7469449606899
$ws.Range("C22").Value = -0.497552364060549
$ws.Range("D22").Value = 11.0817784096294
$ws.Range("B23").Value = 5.31794627982536
$ws.Range("C23").Value = 3.75558626865022
$ws.Range("D23").Value = 6.88030629100049
$ws.Range("C24").Value = 0.272465475975257
$ws.Range("D24").Value = 4.92527904432211
$ws.Range("B25").Value = 7.54256520205381
$ws.Range("C25").Value = 6.39900094752147
$ws.Range("D25").Value = 8.68612945658616
$ws.Range("B26").Value = 17.6477970206699
$ws.Range("C26").Value = 16.5979346389398
$ws.Range("D26").Value = 18.6976594023999
$ws.Range("C27").Value = 13.9455444390322
$ws.Range("D27").Value = 16.9417433936807
$ws.Range("C28").Value = 12.0431941048941
$ws.Range("D28").Value = 19.3572522486888
$ws.Range("C29").Value = 11.3118394310184
$ws.Range("D29").Value = 17.4160035939465
$ws.Range("C30").Value = 10.6248461985365
$ws.Range("D30").Value = 17.5390231245385
$ws.Range("C31").Value = 10.936600861532
$ws.Range("D31").Value = 20.3981277438284
$ws.Range("C32").Value = 9.0947484647813
$ws.Range("D32").Value = 16.9398580953153
$ws.Range("C33").Value = 15.4789182348919
$ws.Range("D33").Value = 23.2654448244822
$ws.Range("C34").Value = 6.18028143959322
$ws.Range("D34").Value = 27.1866627862109
$ws.Range("B35").Value = 21.021707247996
$ws.Range("C35").Value = 18.1301359041404
$ws.Range("D35").Value = 23.9132785918517
$ws.Range("B36").Value = 8.13022396178497
$ws.Range("C36").Value = 3.98836805366789
$ws.Range("D36").Value = 12.2720798699021
$ws.Range("B37").Value = 19.7266556784379
$ws.Range("C37").Value = 17.9362189940129
$ws.Range("D37").Value = 21.5170923628629
$ws.Range("B38").Value = 13.1792237295119
$ws.Range("C38").Value = 11.6935054219436
$ws.Range("D38").Value = 14.6649420370802
$ws.Range("C39").Value = 7.91635295298322
$ws.Range("D39").Value = 11.9013539241769
$ws.Range("C40").Value = 6.39761562820731
$ws.Range("D40").Value = 16.9707838170679
$ws.Range("C41").Value = 5.59252596812335
$ws.Range("D41").Value = 13.907183665388
$ws.Range("C42").Value = 0.793814429595695
$ws.Range("D42").Value = 7.08366746060379
$ws.Range("C43").Value = 4.55848231261023
$ws.Range("D43").Value = 17.3204052331015
$ws.Range("C44").Value = 2.40025769141923
$ws.Range("D44").Value = 11.8774606599113
$ws.Range("C45").Value = 9.82045194890923
$ws.Range("D45").Value = 21.3512160200516
$ws.Range("C46").Value = -4.03241258840341
$ws.Range("D46").Value = 25.4160720547554
$ws.Range("B47").Value = 21.0589466812593
$ws.Range("C47").Value = 16.5187755491723
$ws.Range("D47").Value = 25.5991178133463
$ws.Range("C48").Value = -1.46722063444731
$ws.Range("D48").Value = 4.81273681990517
$ws.Range("B49").Value = 14.669688875493
$ws.Range("C49").Value = 12.139684452525
$ws.Range("D49").Value = 17.1996932984609
$ws.Range("B50").Value = 11.2540732297097
$ws.Range("C50").Value = 9.84501867788117
$ws.Range("D50").Value = 12.6631277815382
$ws.Range("C51").Value = 6.8656709496908
$ws.Range("D51").Value = 10.6562087848265
$ws.Range("C52").Value = 4.60174187825822
$ws.Range("D52").Value = 13.8981237342093
$ws.Range("C53").Value = 4.21321939414173
$ws.Range("D53").Value = 11.5165446565107
$ws.Range("C54").Value = 3.99968172686712
$ws.Range("D54").Value = 13.1158185608224
$ws.Range("C55").Value = 2.65106166928492
$ws.Range("D55").Value = 14.3503907329828
$ws.Range("C56").Value = 0.939417389841765
$ws.Range("D56").Value = 10.2644158132686
$ws.Range("C57").Value = 5.96583135724489
$ws.Range("D57").Value = 16.1136727570489
$ws.Range("C58").Value = -0.633238081160514
$ws.Range("D58").Value = 31.1230279840102
$ws.Range("B59").Value = 14.8639556001856
$ws.Range("C59").Value = 10.900451333537
$ws.Range("D59").Value = 18.8274598668342
$ws.Range("C60").Value = -1.04060259278769
$ws.Range("D60").Value = 8.51652758631026
$ws.Range("B61").Value = 13.3365916496503
$ws.Range("C61").Value = 10.8273751924703
$ws.Range("D61").Value = 15.8458081068302
$ws.Range("B62").Value = 10.3768315010754
$ws.Range("C62").Value = 9.00794967212097
$ws.Range("D62").Value = 11.7457133300299
$ws.Range("C63").Value = 6.96100736420765
$ws.Range("D63").Value = 10.7687426286482
$ws.Range("C64").Value = 8.81468246355078
$ws.Range("D64").Value = 19.8035514288755
$ws.Range("C65").Value = 5.22815903817397
$ws.Range("D65").Value = 13.5492892671897
$ws.Range("C66").Value = -0.184879792697629
$ws.Range("D66").Value = 4.22853913447407
$ws.Range("C67").Value = 1.56094840477654
$ws.Range("D67").Value = 11.5989746194606
$ws.Range("C68").Value = -0.178352945509501
$ws.Range("D68").Value = 6.59865116515342
$ws.Range("C69").Value = 9.1668395528743
$ws.Range("D69").Value = 21.4202731178359
$ws.Range("C70").Value = -2.07734230443387
$ws.Range("D70").Value = 27.1432545114313
$ws.Range("B71").Value = 13.6477580825677
$ws.Range("C71").Value = 9.47829354250712
$ws.Range("D71").Value = 17.8172226226282
$ws.Range("C72").Value = -1.31203552690905
$ws.Range("D72").Value = 4.22490650005778
$ws.Range("B73").Value = 11.892280434608
$ws.Range("C73").Value = 9.52057994586647
$ws.Range("D73").Value = 14.2639809233495
$ws.Range("B74").Value = 11.4867973306881
$ws.Range("C74").Value = 10.0776924518481
$ws.Range("D74").Value = 12.8959022095281
$ws.Range("C75").Value = 6.93554613880903
$ws.Range("D75").Value = 10.6525742034473
$ws.Range("C76").Value = 5.49781564632674
$ws.Range("D76").Value = 15.0373045998489
$ws.Range("C77").Value = 3.57467521743925
$ws.Range("D77").Value = 10.4311744018194
$ws.Range("C78").Value = 0.291442517285037
$ws.Range("D78").Value = 6.54231520834016
$ws.Range("C79").Value = 2.48278517942485
$ws.Range("D79").Value = 14.3884933273422
$ws.Range("C80").Value = 5.90858043963805
$ws.Range("D80").Value = 17.6911370381408
$ws.Range("C81").Value = 6.68526070702708
$ws.Range("D81").Value = 16.2482559869996
$ws.Range("C82").Value = -2.14874222224584
$ws.Range("D82").Value = 27.7878655803572
$ws.Range("B83").Value = 15.6667838528888
$ws.Range("C83").Value = 11.3987964611961
$ws.Range("D83").Value = 19.9347712445815
$ws.Range("C84").Value = -0.547573137192276
$ws.Range("D84").Value = 9.58703195497972
$ws.Range("B85").Value = 13.7977706883562
$ws.Range("C85").Value = 11.2950219042894
$ws.Range("D85").Value = 16.300519472423
